$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; this shifts the existing rows 35-98
# down to 36-99 (matching the target dimension A1:T99).
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new record.
$ws.Range("A35").Value = 1
$ws.Range("B35").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C35").Value = "Arica y Parinacota"
$ws.Range("D35").Value = 44725
$ws.Range("E35").Value = 15
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100106
$ws.Range("H35").Value = "Oleaginosos"
$ws.Range("I35").Value = 100106002
$ws.Range("J35").Value = "Palta"
$ws.Range("K35").Value = "Hass"
$ws.Range("L35").Value = "Segunda"
$ws.Range("M35").Value = 400
$ws.Range("N35").Value = 16000
$ws.Range("O35").Value = 17000
$ws.Range("P35").Value = 16500
$ws.Range("Q35").Value = "$/bandeja 10 kilos"
$ws.Range("R35").Value = "Perú"
$ws.Range("S35").Value = 1650
$ws.Range("T35").Value = 10
